$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Step 1: Remove the "Meta description" paragraph that follows the
#     title heading ("Play 9 Burning Stars for Free - Slot Game Review").
#     It consists of an empty run, a bold "Meta description" run, and a
#     plain text run with the description. Locate it with Find (robust to
#     position), expand to the full paragraph (including its mark), and
#     delete it outright.
$metaFind = $d.Content.Duplicate
$metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$metaFind.Expand(4) | Out-Null
$metaFind.Delete()

# --- Step 2: Replace the closing "For the feature image ..." paragraph
#     with two new paragraphs: a bold title paragraph followed by an
#     italic paragraph holding the (former) meta-description text.
$imgFind = $d.Content.Duplicate
$imgFind.Find.Execute("For the feature image fitting the game", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$imgFind.Expand(4) | Out-Null

$newXml = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 9 Burning Stars for Free - Slot Game Review</w:t></w:r></w:p>" + `
          "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read a review of the slot game 9 Burning Stars, including pros and cons. Try it for free and enjoy the elaborate graphics and engaging music.</w:t></w:r></w:p>"
$imgFind.InsertXML($newXml) | Out-Null
